# Update "想去人数" (F column) figures across sheets, matching the
# upstream data refresh captured in the commit.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value  = 1690
$ws1.Range("F13").Value = 2667
$ws1.Range("F16").Value = 7087
$ws1.Range("F18").Value = 7237
$ws1.Range("F29").Value = 879
$ws1.Range("F31").Value = 286
$ws1.Range("F33").Value = 2430
$ws1.Range("F34").Value = 1204
$ws1.Range("F35").Value = 2736
$ws1.Range("F36").Value = 33
$ws1.Range("F42").Value = 480
$ws1.Range("F43").Value = 528

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F12").Value = 158

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value  = 1690
$ws4.Range("F14").Value = 2667
$ws4.Range("F19").Value = 7087
$ws4.Range("F21").Value = 7237
$ws4.Range("F33").Value = 879
$ws4.Range("F35").Value = 286
$ws4.Range("F37").Value = 2430
$ws4.Range("F38").Value = 1204
$ws4.Range("F40").Value = 2736
$ws4.Range("F41").Value = 33
$ws4.Range("F48").Value = 480
$ws4.Range("F49").Value = 528
